# Auto-update: GitHub Admin Log for remove-repo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 16
$ws.Cells.Item($row, 1).Value = "2025-08-17 16:20:19"
$ws.Cells.Item($row, 2).Value = "remove-repo"
$ws.Cells.Item($row, 3).Value = "new-organization97"
$ws.Cells.Item($row, 4).Value = "Devops"
$ws.Cells.Item($row, 5).Value = "deerepo"

# Column I holds the literal text "False" (not a Boolean). A direct
# Value assignment of the string "False"/"True" is auto-typed to a
# Boolean by Excel, so route it through a formula + paste-values so the
# text result is preserved verbatim.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '="False"'
$scratch.Copy()
$ws.Cells.Item($row, 9).PasteSpecial(-4163)
$scratch.ClearContents()
